$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 284; this pushes the existing rows 284:341 down to 285:342
# and keeps formatting (e.g. the date style on column D) consistent with the
# surrounding rows.
$ws.Rows("284:284").Insert()

# Populate the newly inserted row 284 with the new price-report record.
$ws.Range("A284").Value = 11
$ws.Range("B284").Value = "Vega Monumental Concepción"
$ws.Range("C284").Value = "Bíobío"
$ws.Range("D284").Value = 44694
$ws.Range("E284").Value = 8
$ws.Range("F284").Value = "Fruta"
$ws.Range("G284").Value = 100101
$ws.Range("H284").Value = "Berries"
$ws.Range("I284").Value = 100112025
$ws.Range("J284").Value = "Frutilla"
$ws.Range("K284").Value = "Sin especificar"
$ws.Range("L284").Value = "Primera"
$ws.Range("M284").Value = 160
$ws.Range("N284").Value = 12000
$ws.Range("O284").Value = 13000
$ws.Range("P284").Value = 12500
$ws.Range("Q284").Value = "$/bandeja 7 kilos"
$ws.Range("R284").Value = "Provincia de Melipilla"
$ws.Range("S284").Value = 1786
$ws.Range("T284").Value = 7
